$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("E2").Value = 0.00029023
$ws.Range("F2").Value = 0.020966341
$ws.Range("G2").Value = 0.0004643262269

$ws.Range("E3").Value = 0.004447704
$ws.Range("F3").Value = 0.011329634
$ws.Range("G3").Value = 0.005123980320000001

$ws.Range("E4").Value = 0.007879484
$ws.Range("F4").Value = 0.015375549
$ws.Range("G4").Value = 0.009638699485549133
